$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 31
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
# Row 55
$ws.Range("H55").Value = 158.63637
$ws.Range("J55").Value = 221.5
$ws.Range("L55").Value = 221.5
$ws.Range("N55").Value = -649.5
# Row 137
$ws.Range("H137").Value = 1620.1305
$ws.Range("I137").Value = 1100
$ws.Range("J137").Value = 1803.7059
$ws.Range("K137").Value = 3300
$ws.Range("L137").Value = 5411.1177
$ws.Range("M137").Value = -750
$ws.Range("N137").Value = -10511.1177
# Row 141
$ws.Range("H141").Value = 1905.7142
$ws.Range("I141").Value = 1990.7693
$ws.Range("J141").Value = 800
$ws.Range("K141").Value = 5972.3079
$ws.Range("L141").Value = 2400
$ws.Range("M141").Value = -792.3078999999998
$ws.Range("N141").Value = -12760

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3304.72
$ws.Range("I32").Value = 2649.5144
$ws.Range("K32").Value = 2649.5144
$ws.Range("M32").Value = -2362.5144
# Row 54
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
# Row 61
$ws.Range("H61").Value = 1843.8823
$ws.Range("I61").Value = 1298.7333
$ws.Range("K61").Value = 1298.7333
$ws.Range("M61").Value = -1086.7333
# Row 74
$ws.Range("H74").Value = 1167.1212
$ws.Range("I74").Value = 1138.8182
$ws.Range("J74").Value = 1223.7273
$ws.Range("K74").Value = 1138.8182
$ws.Range("L74").Value = 1223.7273
$ws.Range("M74").Value = -264.8181999999999
$ws.Range("N74").Value = -2971.7273
# Row 77
$ws.Range("H77").Value = 1167.1212
$ws.Range("I77").Value = 1138.8182
$ws.Range("J77").Value = 1223.7273
$ws.Range("K77").Value = 5694.090999999999
$ws.Range("L77").Value = 6118.636500000001
$ws.Range("M77").Value = -1326.090999999999
$ws.Range("N77").Value = -14854.6365
# Row 102
$ws.Range("H102").Value = 2266.6667
$ws.Range("I102").Value = 1800
$ws.Range("K102").Value = 1800
$ws.Range("M102").Value = -178
# Row 110
$ws.Range("H110").Value = 928.95
$ws.Range("I110").Value = 839.17645
$ws.Range("J110").Value = 1437.6666
$ws.Range("K110").Value = 839.17645
$ws.Range("L110").Value = 1437.6666
$ws.Range("M110").Value = 1205.82355
$ws.Range("N110").Value = -5527.6666
# Row 132
$ws.Range("H132").Value = 2240.8865
$ws.Range("I132").Value = 1825.1621
$ws.Range("J132").Value = 4438.2856
$ws.Range("K132").Value = 5475.4863
$ws.Range("L132").Value = 13314.8568
$ws.Range("M132").Value = -2945.4863
$ws.Range("N132").Value = -18374.8568
# Row 136
$ws.Range("H136").Value = 1843.8823
$ws.Range("I136").Value = 1298.7333
$ws.Range("K136").Value = 3896.199900000001
$ws.Range("M136").Value = -1346.199900000001

$ws = $wb.Worksheets.Item("BSM")
# Row 36
$ws.Range("H36").Value = 975
$ws.Range("I36").Value = 975
$ws.Range("K36").Value = 975
$ws.Range("M36").Value = -441

$ws = $wb.Worksheets.Item("CRP")
# Row 23
$ws.Range("H23").Value = 55000
$ws.Range("J23").Value = 60000
$ws.Range("L23").Value = 60000
$ws.Range("N23").Value = -60480
# Row 27
$ws.Range("H27").Value = 55000
$ws.Range("J27").Value = 60000
$ws.Range("L27").Value = 60000
$ws.Range("N27").Value = -60384
# Row 31
$ws.Range("H31").Value = 2054.8845
$ws.Range("I31").Value = 1419.7059
$ws.Range("J31").Value = 3254.6667
$ws.Range("K31").Value = 1419.7059
$ws.Range("L31").Value = 3254.6667
$ws.Range("M31").Value = -1124.7059
$ws.Range("N31").Value = -3844.6667
# Row 34
$ws.Range("H34").Value = 2054.8845
$ws.Range("I34").Value = 1419.7059
$ws.Range("J34").Value = 3254.6667
$ws.Range("K34").Value = 1419.7059
$ws.Range("L34").Value = 3254.6667
$ws.Range("M34").Value = -1217.7059
$ws.Range("N34").Value = -3658.6667

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 26.894737
$ws.Range("I12").Value = 12.090909
$ws.Range("K12").Value = 36.272727
$ws.Range("M12").Value = 136.727273
# Row 68
$ws.Range("H68").Value = 995.1
$ws.Range("I68").Value = 784.53125
$ws.Range("J68").Value = 1369.4445
$ws.Range("K68").Value = 2353.59375
$ws.Range("L68").Value = 4108.333500000001
$ws.Range("M68").Value = -1542.59375
$ws.Range("N68").Value = -5730.333500000001
# Row 71
$ws.Range("H71").Value = 995.1
$ws.Range("I71").Value = 784.53125
$ws.Range("J71").Value = 1369.4445
$ws.Range("K71").Value = 7060.78125
$ws.Range("L71").Value = 12325.0005
$ws.Range("M71").Value = -3004.78125
$ws.Range("N71").Value = -20437.0005
# Row 131
$ws.Range("H131").Value = 3939.465
$ws.Range("I131").Value = 1000
$ws.Range("J131").Value = 4082.8538
$ws.Range("K131").Value = 3000
$ws.Range("L131").Value = 12248.5614
$ws.Range("M131").Value = 2040
$ws.Range("N131").Value = -22328.5614

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6894.1177
$ws.Range("I70").Value = 7600
$ws.Range("K70").Value = 7600
$ws.Range("M70").Value = -7330
# Row 73
$ws.Range("H73").Value = 6894.1177
$ws.Range("I73").Value = 7600
$ws.Range("K73").Value = 7600
$ws.Range("M73").Value = -6664
# Row 102
$ws.Range("H102").Value = 3324.4167
$ws.Range("I102").Value = 996.5
$ws.Range("K102").Value = 996.5
$ws.Range("M102").Value = 625.5
# Row 132
$ws.Range("H132").Value = 2756.1226
$ws.Range("I132").Value = 2170.9285
$ws.Range("J132").Value = 3536.3809
$ws.Range("K132").Value = 6512.7855
$ws.Range("L132").Value = 10609.1427
$ws.Range("M132").Value = -3982.7855
$ws.Range("N132").Value = -15669.1427

$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 2002.3422
$ws.Range("I136").Value = 1361.9656
$ws.Range("J136").Value = 4065.7778
$ws.Range("K136").Value = 4085.8968
$ws.Range("L136").Value = 12197.3334
$ws.Range("M136").Value = -1535.8968
$ws.Range("N136").Value = -17297.3334

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 92179.55
$ws.Range("I126").Value = 112308.336
$ws.Range("J126").Value = 1600
$ws.Range("K126").Value = 336925.008
$ws.Range("L126").Value = 4800
$ws.Range("M126").Value = -334455.008
$ws.Range("N126").Value = -9740
# Row 132
$ws.Range("H132").Value = 15153598
$ws.Range("I132").Value = 23811230
$ws.Range("K132").Value = 71433690
$ws.Range("M132").Value = -71431160
# Row 136
$ws.Range("H136").Value = 10753932
$ws.Range("I136").Value = 16667305
$ws.Range("J136").Value = 2343.9092
$ws.Range("K136").Value = 50001915
$ws.Range("L136").Value = 7031.7276
$ws.Range("M136").Value = -49999365
$ws.Range("N136").Value = -12131.7276

